$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.997.06"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").Value = "3.508.60"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'585.58"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "'132.39"
$ws.Range("E6").Value = "  -1.33%  "
$ws.Range("D7").Value = "3.505.59"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.485"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").Value = "'7.12"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").Value = "4.102.29"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "'27.42"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000178"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.498.85"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "64.009.49"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "'9.79"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("D20").Value = "'13.92"
$ws.Range("E20").Value = "  -3.05%  "
$ws.Range("D21").Value = "'5.60"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'383.10"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").Value = "'0.571"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "3.645.54"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'73.80"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("D29").Value = "'1.57"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.46"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "'8.38"
$ws.Range("E32").Value = "  +0.81%  "
$ws.Range("E33").Value = "  -2.11%  "
$ws.Range("D34").Value = "3.516.96"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'23.50"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("D37").Value = "'0.145"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "'6.89"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'160.33"
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("E42").Value = "  -2.75%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'26.69"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.810"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'41.62"
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("E47").Value = "  -4.13%  "
$ws.Range("D48").Value = "'4.40"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'1.61"
$ws.Range("E49").Value = "  -3.29%  "
$ws.Range("D50").Value = "2.478.84"
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").Value = "'6.78"
$ws.Range("E51").Value = "  -1.77%  "
